$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")

# New score data for columns H:M (rows 12-21); column N = SUM(Hn:Mn)
$data = @{
    12 = @(68, 49, 153, 140, 19.5, 184)
    13 = @(65, 46, 166, 143, 11.5, 171)
    14 = @(61, 43, 143, 136, 2,    141)
    15 = @(60, 56, 153, 146, 18,   116)
    16 = @(54, 66, 163, 137, 13,   161)
    17 = @(59, 58, 185, 156, 12.5, 191)
    18 = @(58, 51, 144, 133, 11.5, 169)
    19 = @(66, 51, 173, 149, 21,   159)
    20 = @(68, 65, 170, 163, 10.5, 216)
    21 = @(52, 53, 131, 116, 21,   161)
}

$cols = @("H", "I", "J", "K", "L", "M")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
    $ws.Range("N$row").Formula = "=SUM(H$row`:M$row)"
}

# Restore the active selection as left by the author
$ws.Range("I19").Select() | Out-Null
